$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing rows down
$ws.Rows.Item(2).Insert()

# New row 2: URL
$ws.Range("A2").Value = "URL"
$ws.Range("B2").Value = "https://preprod-matrix.epixel.link/en/register/"

# Update values that changed after the shift
$ws.Range("B3").Value = "release-mpfp-matrix-business-admin"
$ws.Range("B6").Value = "eyuu0231"
$ws.Range("B7").Value = "abhishna6911@mailinator.com"
$ws.Range("B11").Value = "abhoii8761"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "43436701"
$ws.Range("B17").Value = "enrollment-package-9"
